# Append two new trading-log rows (42, 43) to Sheet1, matching the
# latest ENA trading attempt + resulting opened position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 42: TRADING_ATTEMPT for ENA
$ws.Range("A42").Value = "2025-09-23T06:53:49.541389"
$ws.Range("B42").Value = "TRADING_ATTEMPT"
$ws.Range("C42").Value = "ENA"
$ws.Range("D42").Value = "UNKNOWN"
$ws.Range("E42").Value = 0.6000566049103071
$ws.Range("F42").Value = ""
$ws.Range("G42").Value = ""
$ws.Range("H42").Value = ""
$ws.Range("I42").Value = ""
$ws.Range("J42").Value = ""
$ws.Range("K42").Value = "ATTEMPT"
$ws.Range("L42").Value = "Attempting trade 1/1"

# Row 43: POSITION_OPENED for ENA
$ws.Range("A43").Value = "2025-09-23T06:53:51.349247"
$ws.Range("B43").Value = "POSITION_OPENED"
$ws.Range("C43").Value = "ENA"
$ws.Range("D43").Value = "UNKNOWN"
$ws.Range("E43").Value = 0.6000566049103071
$ws.Range("F43").Value = 2400
$ws.Range("G43").Value = 10
$ws.Range("H43").Value = 0.8506466283814236
$ws.Range("I43").Value = ""
$ws.Range("J43").Value = ""
$ws.Range("K43").Value = "SUCCESS"
$ws.Range("L43").Value = ""
